# Update "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 1161
$sheet1.Range("F3").Value = 415
$sheet1.Range("F4").Value = 264
$sheet1.Range("F5").Value = 146
$sheet1.Range("F6").Value = 9
$sheet1.Range("F7").Value = 12239
$sheet1.Range("F10").Value = 2
$sheet1.Range("F11").Value = 136
$sheet1.Range("F12").Value = 12027
$sheet1.Range("F13").Value = 4806
$sheet1.Range("F14").Value = 2646
$sheet1.Range("F15").Value = 115
$sheet1.Range("F17").Value = 422
$sheet1.Range("F18").Value = 93
$sheet1.Range("F19").Value = 944
$sheet1.Range("F20").Value = 1
$sheet1.Range("F22").Value = 163
$sheet1.Range("F24").Value = 5214

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 1161
$sheet4.Range("F3").Value = 415
$sheet4.Range("F4").Value = 264
$sheet4.Range("F5").Value = 146
$sheet4.Range("F8").Value = 9
$sheet4.Range("F9").Value = 12239
$sheet4.Range("F12").Value = 2
$sheet4.Range("F13").Value = 136
$sheet4.Range("F14").Value = 12027
$sheet4.Range("F15").Value = 4806
$sheet4.Range("F16").Value = 2646
$sheet4.Range("F17").Value = 115
$sheet4.Range("F19").Value = 422
$sheet4.Range("F20").Value = 93
$sheet4.Range("F21").Value = 944
$sheet4.Range("F22").Value = 1
$sheet4.Range("F24").Value = 163
$sheet4.Range("F26").Value = 5214
